# Update "想去人数" (column F) counts across the four sheets of the
# "北京-漫展信息" workbook to reflect newly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$updates = @{
    3  = 24
    4  = 749
    5  = 1328
    6  = 44
    7  = 3211
    9  = 181
    11 = 8253
    12 = 440
    14 = 123
    15 = 66
    17 = 297
    18 = 65
    19 = 331
    20 = 10365
    23 = 37
    25 = 382
    26 = 169
    27 = 144
    28 = 87
    29 = 130
    30 = 2650
    31 = 2065
    32 = 33
    33 = 34
    34 = 2106
    35 = 856
    37 = 265
    38 = 2419
    39 = 2829
    41 = 149
    42 = 315
    43 = 248
    44 = 31
    45 = 95
    46 = 92
    47 = 77
    49 = 58
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# --- Sheet "演出" (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$updates = @{
    10 = 27
    11 = 6
    14 = 30
    19 = 26
    20 = 38
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# --- Sheet "本地生活" (sheet3) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 8

# --- Sheet "全部类型" (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$updates = @{
    4  = 749
    5  = 1328
    6  = 44
    7  = 3211
    9  = 181
    11 = 8253
    12 = 440
    14 = 123
    15 = 66
    17 = 297
    18 = 65
    19 = 331
    20 = 10365
    23 = 37
    25 = 382
    26 = 169
    27 = 30
    28 = 144
    29 = 87
    30 = 130
    31 = 2650
    32 = 2065
    33 = 2106
    34 = 856
    36 = 265
    37 = 2422
    38 = 2830
    39 = 26
    41 = 150
    42 = 315
    43 = 248
    44 = 31
    45 = 95
    46 = 92
    47 = 77
    49 = 58
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
